$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($Range, [string]$Text)
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
}

Set-CellText $ws.Range('D2') '56.642.67'
Set-CellText $ws.Range('E2') '  +0.02%  '
Set-CellText $ws.Range('D3') '2.335.51'
Set-CellText $ws.Range('E3') '  +0.57%  '
Set-CellText $ws.Range('E4') '  -0.10%  '
Set-CellText $ws.Range('D5') '514.35'
Set-CellText $ws.Range('E5') '  -0.41%  '
Set-CellText $ws.Range('D6') '133.35'
Set-CellText $ws.Range('E6') '  -0.12%  '
Set-CellText $ws.Range('E7') '  +0.03%  '
Set-CellText $ws.Range('D8') '0.532'
Set-CellText $ws.Range('E8') '  -0.28%  '
Set-CellText $ws.Range('E9') '  -2.00%  '
Set-CellText $ws.Range('E10') '  -0.84%  '
Set-CellText $ws.Range('D11') '5.33'
Set-CellText $ws.Range('E11') '  +1.42%  '
Set-CellText $ws.Range('E12') '  -0.19%  '
Set-CellText $ws.Range('D13') '23.68'
Set-CellText $ws.Range('E13') '  +0.08%  '
Set-CellText $ws.Range('D14') '2.751.72'
Set-CellText $ws.Range('E14') '  -0.29%  '
Set-CellText $ws.Range('D15') '56.611.30'
Set-CellText $ws.Range('E15') '  +0.13%  '
Set-CellText $ws.Range('E16') '  -0.50%  '
Set-CellText $ws.Range('D17') '2.339.31'
Set-CellText $ws.Range('E17') '  -0.71%  '
Set-CellText $ws.Range('D18') '10.40'
Set-CellText $ws.Range('E18') '  +0.33%  '
Set-CellText $ws.Range('D19') '324.48'
Set-CellText $ws.Range('E19') '  +1.81%  '
Set-CellText $ws.Range('D20') '4.17'
Set-CellText $ws.Range('E20') '  -1.57%  '
Set-CellText $ws.Range('D21') '6.65'
Set-CellText $ws.Range('E21') '  +0.90%  '
Set-CellText $ws.Range('D22') '0.998'
Set-CellText $ws.Range('E22') '  -0.05%  '
Set-CellText $ws.Range('D23') '61.20'
Set-CellText $ws.Range('E23') '  +1.00%  '
Set-CellText $ws.Range('D24') '8.64'
Set-CellText $ws.Range('E24') '  +12.10%  '
Set-CellText $ws.Range('D25') '0.163'
Set-CellText $ws.Range('E25') '  +3.04%  '
Set-CellText $ws.Range('D26') '0.998'
Set-CellText $ws.Range('E26') '  -0.15%  '
Set-CellText $ws.Range('D27') '1.30'
Set-CellText $ws.Range('E27') '  +6.54%  '
Set-CellText $ws.Range('D28') '168.22'
Set-CellText $ws.Range('E28') '  -1.37%  '
Set-CellText $ws.Range('D29') '0.0₃0723'
Set-CellText $ws.Range('E29') '  -1.45%  '
Set-CellText $ws.Range('D30') '1.67'
Set-CellText $ws.Range('E30') '  -0.25%  '
Set-CellText $ws.Range('D31') '6.11'
Set-CellText $ws.Range('E31') '  -1.56%  '
Set-CellText $ws.Range('D32') '18.34'
Set-CellText $ws.Range('E32') '  +0.74%  '
Set-CellText $ws.Range('E33') '  -0.05%  '
Set-CellText $ws.Range('D34') '0.998'
Set-CellText $ws.Range('E34') '  +0.14%  '
Set-CellText $ws.Range('D35') '1.27'
Set-CellText $ws.Range('E35') '  +2.74%  '
Set-CellText $ws.Range('D36') '3.97'
Set-CellText $ws.Range('E36') '  +0.23%  '
Set-CellText $ws.Range('D37') '0.884'
Set-CellText $ws.Range('E37') '  -6.12%  '
Set-CellText $ws.Range('E38') '  +2.53%  '
Set-CellText $ws.Range('D39') '38.44'
Set-CellText $ws.Range('E39') '  +2.61%  '
Set-CellText $ws.Range('D40') '150.03'
Set-CellText $ws.Range('E40') '  +8.80%  '
Set-CellText $ws.Range('D41') '0.376'
Set-CellText $ws.Range('E41') '  -0.69%  '
Set-CellText $ws.Range('D42') '3.58'
Set-CellText $ws.Range('E42') '  +0.79%  '
Set-CellText $ws.Range('D43') '279.84'
Set-CellText $ws.Range('E43') '  +1.76%  '
Set-CellText $ws.Range('D44') '5.09'
Set-CellText $ws.Range('E44') '  +1.21%  '
Set-CellText $ws.Range('D45') '0.0924'
Set-CellText $ws.Range('E45') '  -0.46%  '
Set-CellText $ws.Range('D46') '0.0497'
Set-CellText $ws.Range('E46') '  -0.74%  '
Set-CellText $ws.Range('D47') '0.555'
Set-CellText $ws.Range('E47') '  -0.29%  '
Set-CellText $ws.Range('D48') '18.20'
Set-CellText $ws.Range('E48') '  +6.78%  '
Set-CellText $ws.Range('B49') 'Polygon'
Set-CellText $ws.Range('C49') 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-CellText $ws.Range('D49') '0.380'
Set-CellText $ws.Range('E49') '  +0.63%  '
Set-CellText $ws.Range('B50') 'VeChain'
Set-CellText $ws.Range('C50') 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-CellText $ws.Range('D50') '0.0215'
Set-CellText $ws.Range('E50') '  -0.28%  '
Set-CellText $ws.Range('B51') 'EnergySwap'
Set-CellText $ws.Range('C51') 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-CellText $ws.Range('D51') '17.07'
Set-CellText $ws.Range('E51') '  +2.39%  '
